# update activity file upload
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old test data entirely before laying down the new values
$ws.Cells.Clear()

# Write the new cell values (matches the rebuilt sharedStrings / sheetData)
$ws.Range("A4").Value = "sdfsdfsdf"
$ws.Range("B4").Value = "xcsdfsdf"

$ws.Range("D6").Value = "sdfsdfsd"
$ws.Range("E6").Value = "df"

$ws.Range("B8").Value = "f"
$ws.Range("C8").Value = "sdfsdfsdf"

$ws.Range("B10").Value = "sdfsdfsdf"

$ws.Range("D14").Value = "sdf"

# Move the active selection to D6, matching the saved view state
$ws.Range("D6").Select()
